# This script updates the "cryptos" price-tracker worksheet to reflect a
# refreshed data pull (new Price / Volume(1h) values, plus a handful of
# rows whose relative ranking moved, which re-shuffled the Coin/Link/Price
# rows 26-28 and 43-44).
#
# Several of the "Price" strings look exactly like numbers (e.g. "234.64"),
# but in the source data they are plain text cells (no thousands separator,
# sometimes two dots as in "43.313.03", sometimes very small decimals).
# Setting Range.Value directly on such a string makes Excel silently
# reinterpret it as a number, which would both change the stored value
# (e.g. drop trailing/insignificant formatting) and flip the cell's type.
# To avoid that, for any replacement text that Excel would parse as a
# number we temporarily force the cell to Text format ("@") before writing
# the value, then reset the cell style back to "Normal" so no stray
# number-format/style is left behind (matching the original workbook,
# where these cells carry no special style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($address, $text) {
    $ws.Range($address).Value = $text
}

function Set-NumericLookingTextValue($address, $text) {
    # Force text storage so Excel does not auto-convert the numeric-looking
    # string into a real number, then clear the temporary style so the
    # cell's formatting matches the rest of the untouched sheet.
    $ws.Range($address).NumberFormat = "@"
    $ws.Range($address).Value = $text
    $ws.Range($address).Style = "Normal"
}

Set-TextValue "D2" "43.313.03"
Set-TextValue "E2" "  +0.33%  "
Set-TextValue "D3" "2.358.27"
Set-TextValue "E4" "  -0.25%  "
Set-NumericLookingTextValue "D5" "234.64"
Set-TextValue "E5" "  +1.93%  "
Set-TextValue "E6" "  +2.38%  "
Set-NumericLookingTextValue "D7" "72.48"
Set-TextValue "E7" "  +14.57%  "
Set-TextValue "E8" "  -0.13%  "
Set-TextValue "E9" "  +13.69%  "
Set-TextValue "E10" "  +4.11%  "
Set-NumericLookingTextValue "D11" "27.37"
Set-TextValue "E11" "  -0.34%  "
Set-TextValue "D12" "2.713.56"
Set-TextValue "E12" "  +5.65%  "
Set-TextValue "E13" "  +2.30%  "
Set-NumericLookingTextValue "D14" "16.43"
Set-TextValue "E14" "  +8.44%  "
Set-NumericLookingTextValue "D15" "6.33"
Set-TextValue "E15" "  +5.05%  "
Set-TextValue "E16" "  +5.69%  "
Set-TextValue "D17" "2.358.52"
Set-TextValue "E17" "  +5.62%  "
Set-TextValue "D18" "43.343.65"
Set-TextValue "E18" "  +0.55%  "
Set-NumericLookingTextValue "D19" "0.0000100"
Set-TextValue "E19" "  +4.13%  "
Set-NumericLookingTextValue "D20" "75.26"
Set-TextValue "E20" "  +3.50%  "
Set-NumericLookingTextValue "D21" "6.39"
Set-TextValue "E21" "  +5.48%  "
Set-NumericLookingTextValue "D22" "248.54"
Set-TextValue "E22" "  +1.26%  "
Set-NumericLookingTextValue "D23" "3.88"
Set-TextValue "E23" "  +6.37%  "
Set-TextValue "E24" "  +0.05%  "
Set-NumericLookingTextValue "D25" "2.49"
Set-TextValue "E25" "  +3.15%  "
Set-TextValue "B26" "Cosmos"
Set-TextValue "C26" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-NumericLookingTextValue "D26" "10.15"
Set-TextValue "E26" "  +4.85%  "
Set-TextValue "B27" "EthereumClassic"
Set-TextValue "C27" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-NumericLookingTextValue "D27" "22.62"
Set-TextValue "E27" "  +5.40%  "
Set-TextValue "B28" "Toncoin"
Set-TextValue "C28" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-NumericLookingTextValue "D28" "2.21"
Set-TextValue "E28" "  -3.25%  "
Set-NumericLookingTextValue "D29" "172.87"
Set-TextValue "E29" "  -0.18%  "
Set-TextValue "E30" "  +11.01%  "
Set-NumericLookingTextValue "D31" "0.134"
Set-TextValue "E31" "  +4.09%  "
Set-TextValue "E32" "  +3.52%  "
Set-NumericLookingTextValue "D33" "5.02"
Set-TextValue "E33" "  +2.24%  "
Set-NumericLookingTextValue "D34" "0.0695"
Set-TextValue "E34" "  +3.46%  "
Set-NumericLookingTextValue "D35" "5.07"
Set-TextValue "E35" "  +3.76%  "
Set-NumericLookingTextValue "D36" "3.76"
Set-TextValue "E36" "  +5.20%  "
Set-TextValue "E37" "  +8.20%  "
Set-NumericLookingTextValue "D38" "6.55"
Set-TextValue "E38" "  +4.70%  "
Set-TextValue "E39" "  +3.08%  "
Set-NumericLookingTextValue "D40" "19.25"
Set-TextValue "E40" "  +14.07%  "
Set-TextValue "E41" "  +4.08%  "
Set-TextValue "E42" "  -0.20%  "
Set-TextValue "B43" "Aave"
Set-TextValue "C43" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-NumericLookingTextValue "D43" "99.71"
Set-TextValue "E43" "  +3.83%  "
Set-TextValue "B44" "ARBITRUM"
Set-TextValue "C44" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-NumericLookingTextValue "D44" "1.16"
Set-TextValue "E44" "  +10.44%  "
Set-TextValue "E45" "  +1.84%  "
Set-TextValue "E46" "  +2.50%  "
Set-NumericLookingTextValue "D47" "1.21"
Set-TextValue "E47" "  +2.82%  "
Set-TextValue "D48" "1.446.73"
Set-TextValue "E48" "  +0.61%  "
Set-TextValue "D49" "2.584.89"
Set-TextValue "E49" "  +5.75%  "
Set-TextValue "E51" "  -2.06%  "
